# Updates cryptocurrency price / volume(1h) / hora data for rows 2-51
# on the active worksheet, matching the "Updated symbol list" commit.
# Columns: D = Price, E = Volume(1h), G = Hora
# All three columns are stored as text (inlineStr) in the workbook, so we
# force NumberFormat to Text ("@") before assigning each value to avoid
# Excel auto-converting numeric-looking / percentage-looking strings into
# actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D), new Volume(1h) (E), new Hora (G)
# An empty string means "leave this column unchanged".
$updates = @(
    @(2, "274.85", "-1.28%", "2"),
    @(3, "26.79", "-1.54%", "2"),
    @(4, "4.896", "2.29%", "2"),
    @(5, "0.06347", "1.48%", "2"),
    @(6, "6.934", "1.69%", "2"),
    @(7, "3.345", "2.85%", "2"),
    @(8, "1.297", "36.07%", "2"),
    @(9, "0.8779", "0.77%", "2"),
    @(10, "", "0.77%", "2"),
    @(11, "0.05083", "-3.37%", "2"),
    @(12, "0.07382", "1.97%", "2"),
    @(13, "0.03107", "-0.55%", "2"),
    @(14, "0.09045", "-0.02%", "2"),
    @(15, "0.001568", "1.44%", "2"),
    @(16, "0.0006293", "0.89%", "2"),
    @(17, "0.006011", "0.76%", "2"),
    @(18, "3.469", "-0.01%", "2"),
    @(19, "2.273", "-0.09%", "2"),
    @(20, "0.3166", "", "2"),
    @(21, "0.1325", "1.45%", "2"),
    @(22, "3.900", "1.62%", "2"),
    @(23, "0.04368", "1.57%", "2"),
    @(24, "", "0.43%", "2"),
    @(25, "0.003702", "-12.39%", "2"),
    @(26, "", "0.15%", "2"),
    @(27, "0.0001936", "0.29%", "2"),
    @(28, "", "", "2"),
    @(29, "", "", "2"),
    @(30, "", "", "2"),
    @(31, "", "", "2"),
    @(32, "", "", "2"),
    @(33, "", "", "2"),
    @(34, "", "", "2"),
    @(35, "", "", "2"),
    @(36, "", "", "2"),
    @(37, "", "", "2"),
    @(38, "", "", "2"),
    @(39, "", "", "2"),
    @(40, "0.04081", "1.03%", "2"),
    @(41, "0.006597", "6.62%", "2"),
    @(42, "0.1166", "2.29%", "2"),
    @(43, "0.002128", "0.42%", "2"),
    @(44, "0.01210", "-5.19%", "2"),
    @(45, "0.00005323", "3.97%", "2"),
    @(46, "3.106", "107.85%", "2"),
    @(47, "", "-12.95%", "2"),
    @(48, "", "", "2"),
    @(49, "", "", "2"),
    @(50, "", "", "2"),
    @(51, "", "", "2")
)

foreach ($item in $updates) {
    $r = $item[0]
    $newD = $item[1]
    $newE = $item[2]
    $newG = $item[3]

    if ($newD -ne "") {
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $newD
    }

    if ($newE -ne "") {
        $cellE = $ws.Cells.Item($r, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $newE
    }

    if ($newG -ne "") {
        $cellG = $ws.Cells.Item($r, 7)
        $cellG.NumberFormat = "@"
        $cellG.Value = $newG
    }
}
